# Add a "canonical SMILES" column (D) to the microstate list worksheet.
# Column D holds the canonical (non-isomeric) SMILES - i.e. the existing
# "canonical isomeric SMILES" in column C with the cis/trans bond markers
# ("/" and "\") stripped out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Range("D2").Value = "canonical SMILES"

# Fill in column D for each data row (rows 3-13) by stripping the
# isomeric-bond markers "/" and "\" from the column C SMILES string.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
for ($r = 3; $r -le $lastRow; $r++) {
    $smiles = $ws.Cells.Item($r, 3).Value()
    if ($smiles) {
        $canonical = $smiles.Replace("/", "").Replace("\", "")
        $ws.Cells.Item($r, 4).Value = $canonical
    }
}

# Match the column width used for the new column in the target workbook
# (engine quantizes ColumnWidth to 1/6-character steps, so 36.0 is the
# input value that lands closest to the target stored width of
# 36.85546875 characters).
$ws.Columns.Item(4).ColumnWidth = 36
